# Events.xlsx edit: "fix regions on battle map. Added new ReturnToSpotting territory"
# Inserts a new row (row 72) for event e060 "Reset Round" right after the
# existing e054a "MG Firing - Target Selected" row, pushing the remaining
# rows (old 72-81) down to (73-82).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- Insert a new row at position 72 (shifts rows 72:81 down to 73:82) ---
$ws.Rows("72:72").Insert()

# --- Populate the new row with the e060 "Reset Round" event ---
$ws.Range("A72").Value = "e060"

$resetRoundText = "<Bold>e060 Reset Round</Bold> " + "`n" + `
"<InlineUIContainer><Button Content='r4.77' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   " + "`n" + `
"<LineBreak/><LineBreak/>" + "`n" + `
"Since the round did not end, reset and return back to Spotting Round." + "`n" + `
"<LineBreak/><LineBreak/>" + "`n" + `
"                                            <InlineUIContainer><Image Name='Continue60' Height='100' Width='100'></Image></InlineUIContainer>"

$ws.Range("B72").Value = $resetRoundText

# Match the row height used for the new row in the authored workbook
$ws.Rows("72:72").RowHeight = 90

# --- Restore the view's selected cell to where the author left it ---
$ws.Range("B71").Select()
